$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-19 (column A = element name, column B = type)
$data = @(
    @("demand1", "demand"),
    @("demand2", "demand"),
    @("net1", "net"),
    @("net2", "net"),
    @("pv1", "pv"),
    @("pv2", "pv"),
    @("bat1", "bat"),
    @("bat2", "bat"),
    @("CHP1", "CHP"),
    @("CHP2", "CHP"),
    @("solar_th1", "solar_th"),
    @("solar_th2", "solar_th"),
    @("pvt1", "pvt"),
    @("pvt2", "pvt"),
    @("charging_station1", "charging_station"),
    @("charging_station2", "charging_station"),
    @("heat_pump1", "heat_pump"),
    @("heat_pump2", "heat_pump")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
